# Bug fix for lettergrade method: correct the quiz/assignment score in
# column E for rows 7 and 8 on the "grades" sheet from 0.9 to 1.
# The SUMPRODUCT formula in column H recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1
